$d = $word.ActiveDocument

# Replace all occurrences of "In Progress" with "Done" in the Acceptance
# Criteria status table cells (Sprint 3 status update).
$d.Content.Find.Execute("In Progress", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Done", 2)
